$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.474.93"
$ws.Range("E2").Value = "  -6.97%  "
$ws.Range("D3").Value = "3.741.11"
$ws.Range("E3").Value = "  -6.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.49%  "
$ws.Range("D7").Value = "3.730.17"
$ws.Range("E7").Value = "  -6.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.637"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -6.73%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.722"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.46%  "
$ws.Range("E11").Value = "  -10.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.08"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000302"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -10.34%  "
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").Value = "4.340.53"
$ws.Range("E15").Value = "  -6.01%  "
$ws.Range("D16").Value = "3.748.06"
$ws.Range("E16").Value = "  -5.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.56"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -6.86%  "
$ws.Range("E19").Value = "  -6.87%  "
$ws.Range("D21").Value = "68.443.98"
$ws.Range("E21").Value = "  -6.73%  "
$ws.Range("E22").Value = "  -5.91%  "
$ws.Range("E23").Value = "  -5.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.30"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.86%  "
$ws.Range("E25").Value = "  -7.80%  "
$ws.Range("E26").Value = "  -8.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.87"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.00%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.69"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -8.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.85"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.87%  "
$ws.Range("E34").Value = "  -7.93%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.18"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.83%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "44.33"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "613.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("E38").Value = "  -11.95%  "
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("E43").Value = "  -5.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0449"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.69%  "
$ws.Range("E46").Value = "  +3.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -11.04%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.136"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.33%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.74"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -14.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.16"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.39%  "
$ws.Range("D51").Value = "2.744.66"
$ws.Range("E51").Value = "  -2.19%  "
